$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 427, shifting existing rows 427-474 down to 428-475
$ws.Rows.Item(427).Insert()

# Populate the newly inserted row 427 with the new weekly data point
$ws.Cells.Item(427, 1).Value = 3
$ws.Cells.Item(427, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(427, 3).Value = "Coquimbo"
$ws.Cells.Item(427, 4).Value = 44918
$ws.Cells.Item(427, 5).Value = 5
$ws.Cells.Item(427, 6).Value = 100112040
$ws.Cells.Item(427, 7).Value = "Cilantro"
$ws.Cells.Item(427, 8).Value = "Sin especificar"
$ws.Cells.Item(427, 9).Value = "Primera"
$ws.Cells.Item(427, 10).Value = 230
$ws.Cells.Item(427, 11).Value = 5000
$ws.Cells.Item(427, 12).Value = 5500
$ws.Cells.Item(427, 13).Value = 5261
$ws.Cells.Item(427, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(427, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(427, 16).Value = 1754
$ws.Cells.Item(427, 17).Value = 3
$ws.Cells.Item(427, 18).Value = "Hortaliza"
